# Auto-generated edit script applying the Bahamut_Profits market-data refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific rows
# across multiple item-crafting sheets, matching the scheduled runner's data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 90913656
$ws.Range("I86").Value = 4033.8333
$ws.Range("J86").Value = 200005200
$ws.Range("K86").Value = 4033.8333
$ws.Range("L86").Value = 200005200
$ws.Range("M86").Value = -2910.8333
$ws.Range("N86").Value = -200007446

$ws.Range("H89").Value = 90913656
$ws.Range("I89").Value = 4033.8333
$ws.Range("J89").Value = 200005200
$ws.Range("K89").Value = 20169.1665
$ws.Range("L89").Value = 1000026000
$ws.Range("M89").Value = -14553.1665
$ws.Range("N89").Value = -1000037232

$ws.Range("H98").Value = 5251.25
$ws.Range("I98").Value = 5251.25
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 5251.25
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -3753.25
$ws.Range("N98").Value = $null

$ws.Range("H116").Value = 5362.5
$ws.Range("I116").Value = 5725
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 5725
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -2283
$ws.Range("N116").Value = -11884

$ws.Range("H122").Value = 5251.25
$ws.Range("I122").Value = 5251.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15753.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -13303.75
$ws.Range("N122").Value = $null

$ws.Range("H125").Value = 996.05884
$ws.Range("I125").Value = 999.75
$ws.Range("J125").Value = 987.2
$ws.Range("K125").Value = 8997.75
$ws.Range("L125").Value = 8884.800000000001
$ws.Range("M125").Value = -6537.75
$ws.Range("N125").Value = -13804.8

$ws.Range("H133").Value = 45250
$ws.Range("J133").Value = 45250
$ws.Range("L133").Value = 45250
$ws.Range("N133").Value = -55370

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3331.2144
$ws.Range("J2").Value = 3832.1667
$ws.Range("L2").Value = 3832.1667
$ws.Range("N2").Value = -4058.1667

$ws.Range("H32").Value = 3884.95
$ws.Range("I32").Value = 3884.95
$ws.Range("K32").Value = 3884.95
$ws.Range("M32").Value = -3597.95

$ws.Range("H45").Value = 1041.4584
$ws.Range("I45").Value = 962.1875
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 962.1875
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -585.1875
$ws.Range("N45").Value = -1954

$ws.Range("H74").Value = 1000.5278
$ws.Range("I74").Value = 1006.84375
$ws.Range("J74").Value = 950
$ws.Range("K74").Value = 1006.84375
$ws.Range("L74").Value = 950
$ws.Range("M74").Value = -132.84375
$ws.Range("N74").Value = -2698

$ws.Range("H77").Value = 1000.5278
$ws.Range("I77").Value = 1006.84375
$ws.Range("J77").Value = 950
$ws.Range("K77").Value = 5034.21875
$ws.Range("L77").Value = 4750
$ws.Range("M77").Value = -666.21875
$ws.Range("N77").Value = -13486

$ws.Range("H116").Value = 3331.2144
$ws.Range("J116").Value = 3832.1667
$ws.Range("L116").Value = 3832.1667
$ws.Range("N116").Value = -8420.1667

$ws.Range("H122").Value = 1342.6428
$ws.Range("I122").Value = 1153
$ws.Range("J122").Value = 1418.5
$ws.Range("K122").Value = 3459
$ws.Range("L122").Value = 4255.5
$ws.Range("M122").Value = -1009
$ws.Range("N122").Value = -9155.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3331.2144
$ws.Range("J3").Value = 3832.1667
$ws.Range("L3").Value = 3832.1667
$ws.Range("N3").Value = -4060.1667

$ws.Range("H94").Value = 840.93335
$ws.Range("I94").Value = 998.2174
$ws.Range("J94").Value = 324.14285
$ws.Range("K94").Value = 998.2174
$ws.Range("L94").Value = 324.14285
$ws.Range("M94").Value = -547.2174
$ws.Range("N94").Value = -1226.14285

$ws.Range("H107").Value = 6802.227
$ws.Range("I107").Value = 781.5263
$ws.Range("J107").Value = 44933.332
$ws.Range("K107").Value = 781.5263
$ws.Range("L107").Value = 44933.332
$ws.Range("M107").Value = 1138.4737
$ws.Range("N107").Value = -48773.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3026.8572
$ws.Range("I99").Value = 1867.7646
$ws.Range("K99").Value = 1867.7646
$ws.Range("M99").Value = -369.7646

$ws.Range("H122").Value = 1800
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -10900

$ws.Range("H126").Value = 3026.8572
$ws.Range("I126").Value = 1867.7646
$ws.Range("K126").Value = 5603.293799999999
$ws.Range("M126").Value = -3133.293799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 5915.615
$ws.Range("I46").Value = 980.6
$ws.Range("K46").Value = 2941.8
$ws.Range("M46").Value = -2850.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1119
$ws.Range("I102").Value = 1071.1428
$ws.Range("K102").Value = 1071.1428
$ws.Range("M102").Value = 550.8571999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2621.1667
$ws.Range("I7").Value = 2677.6365
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 2677.6365
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -2565.6365
$ws.Range("N7").Value = -2224

$ws.Range("H16").Value = 5621.3335
$ws.Range("I16").Value = 7849.1875
$ws.Range("J16").Value = 1165.625
$ws.Range("K16").Value = 7849.1875
$ws.Range("L16").Value = 1165.625
$ws.Range("M16").Value = -7679.1875
$ws.Range("N16").Value = -1505.625

$ws.Range("H40").Value = 2453.1765
$ws.Range("I40").Value = 2321.7144
$ws.Range("K40").Value = 2321.7144
$ws.Range("M40").Value = -2185.7144

$ws.Range("H93").Value = 1991.2812
$ws.Range("I93").Value = 1455.4375
$ws.Range("J93").Value = 2527.125
$ws.Range("K93").Value = 1455.4375
$ws.Range("L93").Value = 2527.125
$ws.Range("M93").Value = -207.4375
$ws.Range("N93").Value = -5023.125

$ws.Range("H122").Value = 9048.714
$ws.Range("I122").Value = 14159.625
$ws.Range("K122").Value = 42478.875
$ws.Range("M122").Value = -40028.875

$ws.Range("H126").Value = 2621.1667
$ws.Range("I126").Value = 2677.6365
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 8032.9095
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -5562.9095
$ws.Range("N126").Value = -10940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 721.3
$ws.Range("I107").Value = 619.05884
$ws.Range("J107").Value = 1300.6666
$ws.Range("K107").Value = 1857.17652
$ws.Range("L107").Value = 3901.9998
$ws.Range("M107").Value = 62.82348000000002
$ws.Range("N107").Value = -7741.9998

$ws.Range("H113").Value = 299.75
$ws.Range("I113").Value = 256.22223
$ws.Range("J113").Value = 355.7143
$ws.Range("K113").Value = 768.66669
$ws.Range("L113").Value = 1067.1429
$ws.Range("M113").Value = 1401.33331
$ws.Range("N113").Value = -5407.1429

$ws.Range("H131").Value = 44562.223
$ws.Range("J131").Value = 44562.223
$ws.Range("L131").Value = 44562.223
$ws.Range("N131").Value = -54642.223

$ws.Range("H132").Value = 957.06525
$ws.Range("I132").Value = 828.91174
$ws.Range("J132").Value = 1320.1666
$ws.Range("K132").Value = 2486.73522
$ws.Range("L132").Value = 3960.4998
$ws.Range("M132").Value = 43.26477999999997
$ws.Range("N132").Value = -9020.4998
